# Apply scheduled market-data refresh to the Leve profit sheets.
# Values reconstructed from the canonical-OOXML diff (per sheet/row/column).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 80
$ws.Cells.Item(80, 8).Value = 439.18182
$ws.Cells.Item(80, 9).Value = 417.125
$ws.Cells.Item(80, 10).Value = 498
$ws.Cells.Item(80, 11).Value = 1251.375
$ws.Cells.Item(80, 12).Value = 1494
$ws.Cells.Item(80, 13).Value = -253.375
$ws.Cells.Item(80, 14).Value = -3490

# ALC row 83
$ws.Cells.Item(83, 8).Value = 439.18182
$ws.Cells.Item(83, 9).Value = 417.125
$ws.Cells.Item(83, 10).Value = 498
$ws.Cells.Item(83, 11).Value = 3754.125
$ws.Cells.Item(83, 12).Value = 4482
$ws.Cells.Item(83, 13).Value = 1237.875
$ws.Cells.Item(83, 14).Value = -14466

# ALC row 113
$ws.Cells.Item(113, 8).Value = 4409.6
$ws.Cells.Item(113, 9).Value = 4409.6
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = 4409.6
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 13).Value = -1155.6

# ALC row 115
$ws.Cells.Item(115, 8).Value = 380
$ws.Cells.Item(115, 9).Value = 380
$ws.Cells.Item(115, 10).Value = 0
$ws.Cells.Item(115, 11).Value = 1140
$ws.Cells.Item(115, 12).Value = 0
$ws.Cells.Item(115, 13).Value = 427

# ALC row 131
$ws.Cells.Item(131, 8).Value = 7181.25
$ws.Cells.Item(131, 9).Value = 2908.3333
$ws.Cells.Item(131, 10).Value = 20000
$ws.Cells.Item(131, 11).Value = 8724.999899999999
$ws.Cells.Item(131, 12).Value = 60000
$ws.Cells.Item(131, 13).Value = -3684.999899999999
$ws.Cells.Item(131, 14).Value = -70080

$ws = $wb.Worksheets.Item("ARM")
# ARM row 5
$ws.Cells.Item(5, 8).Value = 106.5
$ws.Cells.Item(5, 9).Value = 115.28571
$ws.Cells.Item(5, 10).Value = 86
$ws.Cells.Item(5, 11).Value = 115.28571
$ws.Cells.Item(5, 12).Value = 86
$ws.Cells.Item(5, 13).Value = -3.285709999999995
$ws.Cells.Item(5, 14).Value = -310

# ARM row 61
$ws.Cells.Item(61, 8).Value = 6570.5713
$ws.Cells.Item(61, 9).Value = 5598.8
$ws.Cells.Item(61, 10).Value = 9000
$ws.Cells.Item(61, 11).Value = 5598.8
$ws.Cells.Item(61, 12).Value = 9000
$ws.Cells.Item(61, 13).Value = -5386.8
$ws.Cells.Item(61, 14).Value = -9424

# ARM row 113
$ws.Cells.Item(113, 8).Value = 0
$ws.Cells.Item(113, 9).Value = 0
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = 0
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 14).ClearContents()

# ARM row 122
$ws.Cells.Item(122, 8).Value = 2499.6667
$ws.Cells.Item(122, 9).Value = 2499.6667
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 7499.000100000001
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -5049.000100000001

# ARM row 132
$ws.Cells.Item(132, 8).Value = 1308.7142
$ws.Cells.Item(132, 9).Value = 1447.75
$ws.Cells.Item(132, 10).Value = 474.5
$ws.Cells.Item(132, 11).Value = 4343.25
$ws.Cells.Item(132, 12).Value = 1423.5
$ws.Cells.Item(132, 13).Value = -1813.25
$ws.Cells.Item(132, 14).Value = -6483.5

# ARM row 136
$ws.Cells.Item(136, 8).Value = 6570.5713
$ws.Cells.Item(136, 9).Value = 5598.8
$ws.Cells.Item(136, 10).Value = 9000
$ws.Cells.Item(136, 11).Value = 16796.4
$ws.Cells.Item(136, 12).Value = 27000
$ws.Cells.Item(136, 13).Value = -14246.4
$ws.Cells.Item(136, 14).Value = -32100

$ws = $wb.Worksheets.Item("BSM")
# BSM row 4
$ws.Cells.Item(4, 8).Value = 106.5
$ws.Cells.Item(4, 9).Value = 115.28571
$ws.Cells.Item(4, 10).Value = 86
$ws.Cells.Item(4, 11).Value = 115.28571
$ws.Cells.Item(4, 12).Value = 86
$ws.Cells.Item(4, 13).Value = -0.2857099999999946
$ws.Cells.Item(4, 14).Value = -316

# BSM row 22
$ws.Cells.Item(22, 8).Value = 3863.625
$ws.Cells.Item(22, 9).Value = 4058.4285
$ws.Cells.Item(22, 10).Value = 2500
$ws.Cells.Item(22, 11).Value = 4058.4285
$ws.Cells.Item(22, 12).Value = 2500
$ws.Cells.Item(22, 13).Value = -3885.4285
$ws.Cells.Item(22, 14).Value = -2846

# BSM row 105
$ws.Cells.Item(105, 8).Value = 1133.05
$ws.Cells.Item(105, 9).Value = 1035.3529
$ws.Cells.Item(105, 10).Value = 1686.6666
$ws.Cells.Item(105, 11).Value = 1035.3529
$ws.Cells.Item(105, 12).Value = 1686.6666
$ws.Cells.Item(105, 13).Value = 711.6470999999999
$ws.Cells.Item(105, 14).Value = -5180.6666

$ws = $wb.Worksheets.Item("CRP")
# CRP row 7
$ws.Cells.Item(7, 8).Value = 4150.68
$ws.Cells.Item(7, 9).Value = 5984.8237
$ws.Cells.Item(7, 10).Value = 253.125
$ws.Cells.Item(7, 11).Value = 5984.8237
$ws.Cells.Item(7, 12).Value = 253.125
$ws.Cells.Item(7, 13).Value = -5871.8237
$ws.Cells.Item(7, 14).Value = -479.125

# CRP row 19
$ws.Cells.Item(19, 8).Value = 115
$ws.Cells.Item(19, 9).Value = 115
$ws.Cells.Item(19, 10).Value = 0
$ws.Cells.Item(19, 11).Value = 115
$ws.Cells.Item(19, 12).Value = 0
$ws.Cells.Item(19, 13).Value = 55
$ws.Cells.Item(19, 14).ClearContents()

# CRP row 24
$ws.Cells.Item(24, 8).Value = 115
$ws.Cells.Item(24, 9).Value = 115
$ws.Cells.Item(24, 10).Value = 0
$ws.Cells.Item(24, 11).Value = 115
$ws.Cells.Item(24, 12).Value = 0
$ws.Cells.Item(24, 13).Value = 55
$ws.Cells.Item(24, 14).ClearContents()

# CRP row 51
$ws.Cells.Item(51, 8).Value = 28713.75
$ws.Cells.Item(51, 9).Value = 5363.3335
$ws.Cells.Item(51, 10).Value = 98765
$ws.Cells.Item(51, 11).Value = 5363.3335
$ws.Cells.Item(51, 12).Value = 98765
$ws.Cells.Item(51, 13).Value = -4627.3335
$ws.Cells.Item(51, 14).Value = -100237

# CRP row 61
$ws.Cells.Item(61, 8).Value = 28713.75
$ws.Cells.Item(61, 9).Value = 5363.3335
$ws.Cells.Item(61, 10).Value = 98765
$ws.Cells.Item(61, 11).Value = 5363.3335
$ws.Cells.Item(61, 12).Value = 98765
$ws.Cells.Item(61, 13).Value = -5015.3335
$ws.Cells.Item(61, 14).Value = -99461

# CRP row 62
$ws.Cells.Item(62, 8).Value = 8668
$ws.Cells.Item(62, 9).Value = 8000
$ws.Cells.Item(62, 10).Value = 10004
$ws.Cells.Item(62, 11).Value = 8000
$ws.Cells.Item(62, 12).Value = 10004
$ws.Cells.Item(62, 13).Value = -7376
$ws.Cells.Item(62, 14).Value = -11252

# CRP row 65
$ws.Cells.Item(65, 8).Value = 8668
$ws.Cells.Item(65, 9).Value = 8000
$ws.Cells.Item(65, 10).Value = 10004
$ws.Cells.Item(65, 11).Value = 40000
$ws.Cells.Item(65, 12).Value = 50020
$ws.Cells.Item(65, 13).Value = -36880
$ws.Cells.Item(65, 14).Value = -56260

# CRP row 132
$ws.Cells.Item(132, 8).Value = 2597
$ws.Cells.Item(132, 9).Value = 2597
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 7791
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -5261
$ws.Cells.Item(132, 14).ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# CUL row 113
$ws.Cells.Item(113, 8).Value = 1218.6471
$ws.Cells.Item(113, 9).Value = 986.4286
$ws.Cells.Item(113, 10).Value = 1381.2
$ws.Cells.Item(113, 11).Value = 2959.2858
$ws.Cells.Item(113, 12).Value = 4143.6
$ws.Cells.Item(113, 13).Value = -789.2857999999997
$ws.Cells.Item(113, 14).Value = -8483.6

# CUL row 117
$ws.Cells.Item(117, 8).Value = 695
$ws.Cells.Item(117, 9).Value = 0
$ws.Cells.Item(117, 10).Value = 695
$ws.Cells.Item(117, 11).Value = 0
$ws.Cells.Item(117, 12).Value = 2085
$ws.Cells.Item(117, 14).Value = -8969

$ws = $wb.Worksheets.Item("GSM")
# GSM row 107
$ws.Cells.Item(107, 8).Value = 1927.4
$ws.Cells.Item(107, 9).Value = 1546.7273
$ws.Cells.Item(107, 10).Value = 2974.25
$ws.Cells.Item(107, 11).Value = 1546.7273
$ws.Cells.Item(107, 12).Value = 2974.25
$ws.Cells.Item(107, 13).Value = 373.2727
$ws.Cells.Item(107, 14).Value = -6814.25

$ws = $wb.Worksheets.Item("LTW")
# LTW row 22
$ws.Cells.Item(22, 8).Value = 834
$ws.Cells.Item(22, 9).Value = 810.5714
$ws.Cells.Item(22, 10).Value = 875
$ws.Cells.Item(22, 11).Value = 810.5714
$ws.Cells.Item(22, 12).Value = 875
$ws.Cells.Item(22, 13).Value = -515.5714
$ws.Cells.Item(22, 14).Value = -1465

# LTW row 27
$ws.Cells.Item(27, 8).Value = 834
$ws.Cells.Item(27, 9).Value = 810.5714
$ws.Cells.Item(27, 10).Value = 875
$ws.Cells.Item(27, 11).Value = 810.5714
$ws.Cells.Item(27, 12).Value = 875
$ws.Cells.Item(27, 13).Value = -703.5714
$ws.Cells.Item(27, 14).Value = -1089

# LTW row 61
$ws.Cells.Item(61, 8).Value = 7157.143
$ws.Cells.Item(61, 9).Value = 6033.6665
$ws.Cells.Item(61, 10).Value = 7999.75
$ws.Cells.Item(61, 11).Value = 6033.6665
$ws.Cells.Item(61, 12).Value = 7999.75
$ws.Cells.Item(61, 13).Value = -5831.6665
$ws.Cells.Item(61, 14).Value = -8403.75

# LTW row 96
$ws.Cells.Item(96, 8).Value = 70000
$ws.Cells.Item(96, 9).Value = 0
$ws.Cells.Item(96, 10).Value = 70000
$ws.Cells.Item(96, 11).Value = 0
$ws.Cells.Item(96, 12).Value = 70000
$ws.Cells.Item(96, 14).Value = -75492

# LTW row 105
$ws.Cells.Item(105, 8).Value = 16000
$ws.Cells.Item(105, 9).Value = 0
$ws.Cells.Item(105, 10).Value = 16000
$ws.Cells.Item(105, 11).Value = 0
$ws.Cells.Item(105, 12).Value = 16000
$ws.Cells.Item(105, 14).Value = -22988

# LTW row 113
$ws.Cells.Item(113, 8).Value = 7157.143
$ws.Cells.Item(113, 9).Value = 6033.6665
$ws.Cells.Item(113, 10).Value = 7999.75
$ws.Cells.Item(113, 11).Value = 6033.6665
$ws.Cells.Item(113, 12).Value = 7999.75
$ws.Cells.Item(113, 13).Value = -3863.6665
$ws.Cells.Item(113, 14).Value = -12339.75

# LTW row 132
$ws.Cells.Item(132, 8).Value = 5750
$ws.Cells.Item(132, 9).Value = 5500
$ws.Cells.Item(132, 10).Value = 6000
$ws.Cells.Item(132, 11).Value = 16500
$ws.Cells.Item(132, 12).Value = 18000
$ws.Cells.Item(132, 13).Value = -13970
$ws.Cells.Item(132, 14).Value = -23060

$ws = $wb.Worksheets.Item("WVR")
# WVR row 113
$ws.Cells.Item(113, 8).Value = 417.5
$ws.Cells.Item(113, 9).Value = 223.75
$ws.Cells.Item(113, 10).Value = 805
$ws.Cells.Item(113, 11).Value = 671.25
$ws.Cells.Item(113, 12).Value = 2415
$ws.Cells.Item(113, 13).Value = 1498.75
$ws.Cells.Item(113, 14).Value = -6755

# WVR row 132
$ws.Cells.Item(132, 8).Value = 850.5
$ws.Cells.Item(132, 9).Value = 850.5
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 2551.5
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -21.5

# WVR row 136
$ws.Cells.Item(136, 8).Value = 4061.111
$ws.Cells.Item(136, 9).Value = 2983.3333
$ws.Cells.Item(136, 10).Value = 5138.8887
$ws.Cells.Item(136, 11).Value = 8949.999899999999
$ws.Cells.Item(136, 12).Value = 15416.6661
$ws.Cells.Item(136, 13).Value = -6399.999899999999
$ws.Cells.Item(136, 14).Value = -20516.6661

Write-Output "Applied scheduled-runner price refresh to $($wb.Worksheets.Count) sheets."
